# Travis County 2021 bg SVI - refreshed factor-analysis run:
# the underlying variable ordering / re-fit shifts the significant-component
# groupings and loading-factor figures (and downstream variance summaries)
# by small numeric deltas. Write the refreshed values cell-by-cell.
$wb = $excel.ActiveWorkbook

# --- Significant Components ---
$ws = $wb.Worksheets.Item('Significant Components')
$ws.Range("C2").Value = '[''QEXTRCT'' ''QESL'' ''QHISPC'' ''QEDLESHI'' ''QNOHLTH'' ''PERCAP'' ''QFHH'' ''PPUNIT'']'
$ws.Range("C3").Value = '[''PERCAP'' ''QRICH'' ''MDHSEVAL'']'
$ws.Range("C4").Value = '[''QAGEDEP'' ''MEDAGE'' ''QSSBEN'']'
$ws.Range("C5").Value = '[''QAGEDEP'' ''QFEMLBR'' ''QFEMALE'']'

# --- Loading Factors ---
$ws = $wb.Worksheets.Item('Loading Factors')
$ws.Range("A2").Value = 'QEXTRCT'
$ws.Range("B2").Value = 0.6861621369840734
$ws.Range("C2").Value = 0.1369820409115101
$ws.Range("D2").Value = 0.03705268242987249
$ws.Range("E2").Value = -0.1330634971136807
$ws.Range("F2").Value = 0.02209505719912338
$ws.Range("A3").Value = 'QESL'
$ws.Range("B3").Value = 0.7610062595033303
$ws.Range("C3").Value = 0.1333403856157095
$ws.Range("D3").Value = -0.0457860221660163
$ws.Range("E3").Value = -0.08619628694586934
$ws.Range("F3").Value = 0.1477259463111928
$ws.Range("B4").Value = 0.7657118309845242
$ws.Range("C4").Value = 0.3426616927426591
$ws.Range("D4").Value = -0.1601160553694283
$ws.Range("E4").Value = -0.03265425270417702
$ws.Range("F4").Value = 0.02972290448546306
$ws.Range("A5").Value = 'QEDLESHI'
$ws.Range("B5").Value = 0.8518346711795487
$ws.Range("C5").Value = 0.1560533323342601
$ws.Range("D5").Value = 0.01764598775820661
$ws.Range("E5").Value = -0.04447833179459585
$ws.Range("F5").Value = 0.1021079176823914
$ws.Range("A6").Value = 'QNOHLTH'
$ws.Range("B6").Value = 0.6226791937619538
$ws.Range("C6").Value = 0.3646334545006599
$ws.Range("D6").Value = -0.09045423546840405
$ws.Range("E6").Value = -0.07925175329250797
$ws.Range("F6").Value = 0.2337813356936701
$ws.Range("B7").Value = 0.4581655783272837
$ws.Range("C7").Value = 0.6946266546328257
$ws.Range("D7").Value = -0.1622585102102716
$ws.Range("E7").Value = 0.08440055001674149
$ws.Range("F7").Value = 0.1690618353542978
$ws.Range("B8").Value = 0.4629594412181414
$ws.Range("C8").Value = 0.2259451969006321
$ws.Range("D8").Value = -0.02445972851334618
$ws.Range("E8").Value = 0.266142804994299
$ws.Range("F8").Value = -0.02643591811708856
$ws.Range("B9").Value = 0.5405794971023409
$ws.Range("C9").Value = -0.01941550115363555
$ws.Range("D9").Value = -0.071675345993517
$ws.Range("E9").Value = 0.06034178759289222
$ws.Range("F9").Value = -0.509311952649019
$ws.Range("A10").Value = 'QRICH'
$ws.Range("B10").Value = 0.2158065126519416
$ws.Range("C10").Value = 0.7856952546196763
$ws.Range("D10").Value = -0.1355462535855882
$ws.Range("E10").Value = 0.006122161864189169
$ws.Range("F10").Value = 0.3941286348942635
$ws.Range("A11").Value = 'MDHSEVAL'
$ws.Range("B11").Value = 0.3520833576086306
$ws.Range("C11").Value = 0.7473204774786562
$ws.Range("D11").Value = -0.05970158999541795
$ws.Range("E11").Value = 0.04322327908081265
$ws.Range("F11").Value = -0.02557815698885375
$ws.Range("B12").Value = 0.00008804567514617852
$ws.Range("C12").Value = 0.1850510056404087
$ws.Range("D12").Value = -0.4078865228918835
$ws.Range("E12").Value = -0.05843007178100684
$ws.Range("F12").Value = 0.7697888638425988
$ws.Range("B13").Value = 0.08296100993304176
$ws.Range("C13").Value = 0.03898916580297255
$ws.Range("D13").Value = -0.03483637831037432
$ws.Range("E13").Value = -0.01700687869317976
$ws.Range("F13").Value = 0.605016339611047
$ws.Range("B14").Value = 0.2551100558030827
$ws.Range("C14").Value = 0.157057029416588
$ws.Range("D14").Value = -0.2008460452094055
$ws.Range("E14").Value = 0.0732438610693109
$ws.Range("F14").Value = 0.5089186658371532
$ws.Range("A15").Value = 'QAGEDEP'
$ws.Range("B15").Value = 0.00832034255522143
$ws.Range("C15").Value = -0.05310913538703931
$ws.Range("D15").Value = 0.679832305099877
$ws.Range("E15").Value = 0.6149755400541009
$ws.Range("F15").Value = -0.1056721484969751
$ws.Range("A16").Value = 'MEDAGE'
$ws.Range("B16").Value = -0.2157680421423065
$ws.Range("C16").Value = -0.2398741567895207
$ws.Range("D16").Value = 0.7744618834270887
$ws.Range("E16").Value = -0.07605964727241939
$ws.Range("F16").Value = -0.2120411694857769
$ws.Range("B17").Value = 0.0348349409225177
$ws.Range("C17").Value = -0.02143099166194162
$ws.Range("D17").Value = 0.8056962844545668
$ws.Range("E17").Value = 0.01169225852141011
$ws.Range("F17").Value = -0.1333783606423558
$ws.Range("A18").Value = 'QFEMLBR'
$ws.Range("B18").Value = -0.1508984804554668
$ws.Range("C18").Value = 0.06421827975654623
$ws.Range("D18").Value = -0.0891182476464424
$ws.Range("E18").Value = 0.7138512624185095
$ws.Range("F18").Value = 0.005216870371412786
$ws.Range("A19").Value = 'QFEMALE'
$ws.Range("B19").Value = 0.01407163346046028
$ws.Range("C19").Value = 0.003044396636676153
$ws.Range("D19").Value = 0.1083432321015054
$ws.Range("E19").Value = 0.9214417375769425
$ws.Range("F19").Value = 0.002500953311278286

# --- All Refactor Variances ---
$ws = $wb.Worksheets.Item('All Refactor Variances')
$ws.Range("J2").Value = 4.083072909085918
$ws.Range("K2").Value = 2.086224409112965
$ws.Range("L2").Value = 2.080481394546425
$ws.Range("M2").Value = 2.023412423031613
$ws.Range("N2").Value = 1.866365268053231
$ws.Range("O2").Value = 3.7793943633677
$ws.Range("P2").Value = 2.146393959314242
$ws.Range("Q2").Value = 2.030314401688588
$ws.Range("R2").Value = 1.869786631028111
$ws.Range("S2").Value = 1.824514580245487
$ws.Range("J3").Value = 0.2041536454542959
$ws.Range("K3").Value = 0.1043112204556483
$ws.Range("L3").Value = 0.1040240697273213
$ws.Range("M3").Value = 0.1011706211515806
$ws.Range("N3").Value = 0.09331826340266156
$ws.Range("O3").Value = 0.2099663535204278
$ws.Range("P3").Value = 0.1192441088507912
$ws.Range("Q3").Value = 0.1127952445382549
$ws.Range("R3").Value = 0.1038770350571173
$ws.Range("S3").Value = 0.1013619211247493
$ws.Range("J4").Value = 0.2041536454542959
$ws.Range("K4").Value = 0.3084648659099442
$ws.Range("L4").Value = 0.4124889356372654
$ws.Range("M4").Value = 0.513659556788846
$ws.Range("N4").Value = 0.6069778201915076
$ws.Range("O4").Value = 0.2099663535204278
$ws.Range("P4").Value = 0.329210462371219
$ws.Range("Q4").Value = 0.4420057069094739
$ws.Range("R4").Value = 0.5458827419665913
$ws.Range("S4").Value = 0.6472446630913405
$ws.Range("J5").Value = 0.336344490132907
$ws.Range("K5").Value = 0.1718534302006901
$ws.Range("L5").Value = 0.1713803474639331
$ws.Range("M5").Value = 0.1666792719372518
$ws.Range("N5").Value = 0.1537424602652181
$ws.Range("O5").Value = 0.3244002855390047
$ws.Range("P5").Value = 0.1842334369838802
$ws.Range("Q5").Value = 0.1742698719206542
$ws.Range("R5").Value = 0.1604911418828616
$ws.Range("S5").Value = 0.1566052636735992

# --- Final Variances ---
$ws = $wb.Worksheets.Item('Final Variances')
$ws.Range("B2").Value = 3.7793943633677
$ws.Range("C2").Value = 2.146393959314242
$ws.Range("D2").Value = 2.030314401688588
$ws.Range("E2").Value = 1.869786631028111
$ws.Range("F2").Value = 1.824514580245487
$ws.Range("B3").Value = 0.2099663535204278
$ws.Range("C3").Value = 0.1192441088507912
$ws.Range("D3").Value = 0.1127952445382549
$ws.Range("E3").Value = 0.1038770350571173
$ws.Range("F3").Value = 0.1013619211247493
$ws.Range("B4").Value = 0.2099663535204278
$ws.Range("C4").Value = 0.329210462371219
$ws.Range("D4").Value = 0.4420057069094739
$ws.Range("E4").Value = 0.5458827419665913
$ws.Range("F4").Value = 0.6472446630913405
$ws.Range("B5").Value = 0.3244002855390047
$ws.Range("C5").Value = 0.1842334369838802
$ws.Range("D5").Value = 0.1742698719206542
$ws.Range("E5").Value = 0.1604911418828616
$ws.Range("F5").Value = 0.1566052636735992

# --- Included and Excluded ---
$ws = $wb.Worksheets.Item('Included and Excluded')
$ws.Range("B2").Value = '[[''QEXTRCT'', ''QESL'', ''QHISPC'', ''QEDLESHI'', ''QNOHLTH'', ''PERCAP'', ''QFHH'', ''PPUNIT'', ''QRICH'', ''MDHSEVAL'', ''QAGEDEP'', ''MEDAGE'', ''QSSBEN'', ''QFEMLBR'', ''QFEMALE'', ''QRENTER'', ''QNOAUTO'', ''QPOVTY'']]'

